# Scheduled runner update: refresh computed market/profit figures across
# several sheets (ALC, ARM, BSM, CRP, CUL, LTW).

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# ALC
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H113").Value = 100004220
$ws.Range("I113").Value = 166669380
$ws.Range("J113").Value = 6487.5
$ws.Range("K113").Value = 166669380
$ws.Range("L113").Value = 6487.5
$ws.Range("M113").Value = -166666126
$ws.Range("N113").Value = -12995.5

$ws.Range("H137").Value = 2105.2222
$ws.Range("I137").Value = 1819.4
$ws.Range("J137").Value = 2462.5
$ws.Range("K137").Value = 5458.200000000001
$ws.Range("L137").Value = 7387.5
$ws.Range("M137").Value = -2908.200000000001
$ws.Range("N137").Value = -12487.5

$ws.Range("H138").Value = 1353.5
$ws.Range("I138").Value = 531.1070999999999
$ws.Range("K138").Value = 1593.3213
$ws.Range("M138").Value = 3546.6787

# ----------------------------------------------------------------------
# ARM
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H2").Value = 1517.7307
$ws.Range("I2").Value = 1137.75
$ws.Range("J2").Value = 2125.7
$ws.Range("K2").Value = 1137.75
$ws.Range("L2").Value = 2125.7
$ws.Range("M2").Value = -1024.75
$ws.Range("N2").Value = -2351.7

$ws.Range("H6").Value = 25002500
$ws.Range("I6").Value = 25002500
$ws.Range("K6").Value = 25002500
$ws.Range("M6").Value = -25002327

$ws.Range("H32").Value = 2585.9473
$ws.Range("I32").Value = 1970.6364
$ws.Range("K32").Value = 1970.6364
$ws.Range("M32").Value = -1683.6364

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H102").Value = 1228.1765
$ws.Range("I102").Value = 1062.8572
$ws.Range("J102").Value = 1999.6666
$ws.Range("K102").Value = 1062.8572
$ws.Range("L102").Value = 1999.6666
$ws.Range("M102").Value = 559.1428000000001
$ws.Range("N102").Value = -5243.6666

$ws.Range("H103").Value = 46787.332
$ws.Range("J103").Value = 46787.332
$ws.Range("L103").Value = 46787.332
$ws.Range("N103").Value = -49131.332

$ws.Range("H116").Value = 1517.7307
$ws.Range("I116").Value = 1137.75
$ws.Range("J116").Value = 2125.7
$ws.Range("K116").Value = 1137.75
$ws.Range("L116").Value = 2125.7
$ws.Range("M116").Value = 1156.25
$ws.Range("N116").Value = -6713.7

$ws.Range("H122").Value = 2699.9092
$ws.Range("I122").Value = 2719.9
$ws.Range("K122").Value = 8159.700000000001
$ws.Range("M122").Value = -5709.700000000001

# ----------------------------------------------------------------------
# BSM
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H3").Value = 1517.7307
$ws.Range("I3").Value = 1137.75
$ws.Range("J3").Value = 2125.7
$ws.Range("K3").Value = 1137.75
$ws.Range("L3").Value = 2125.7
$ws.Range("M3").Value = -1023.75
$ws.Range("N3").Value = -2353.7

$ws.Range("H94").Value = 2537.973
$ws.Range("I94").Value = 1173.8889
$ws.Range("J94").Value = 6221
$ws.Range("K94").Value = 1173.8889
$ws.Range("L94").Value = 6221
$ws.Range("M94").Value = -722.8888999999999
$ws.Range("N94").Value = -7123

# ----------------------------------------------------------------------
# CRP
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 16812.5
$ws.Range("I31").Value = 17958.334
$ws.Range("J31").Value = 6500
$ws.Range("K31").Value = 17958.334
$ws.Range("L31").Value = 6500
$ws.Range("M31").Value = -17663.334
$ws.Range("N31").Value = -7090

$ws.Range("H34").Value = 16812.5
$ws.Range("I34").Value = 17958.334
$ws.Range("J34").Value = 6500
$ws.Range("K34").Value = 17958.334
$ws.Range("L34").Value = 6500
$ws.Range("M34").Value = -17756.334
$ws.Range("N34").Value = -6904

$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = 0

$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = 0

$ws.Range("H134").Value = 1047.1818
$ws.Range("I134").Value = 871.8946999999999
$ws.Range("K134").Value = 2615.6841
$ws.Range("M134").Value = -80.68409999999994

# ----------------------------------------------------------------------
# CUL
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H4").Value = 74.7
$ws.Range("I4").Value = 74.7
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 224.1
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -112.1

$ws.Range("H22").Value = 6294.4443
$ws.Range("I22").Value = 9727.272000000001
$ws.Range("K22").Value = 29181.816
$ws.Range("M22").Value = -29012.816

$ws.Range("H27").Value = 6294.4443
$ws.Range("I27").Value = 9727.272000000001
$ws.Range("K27").Value = 29181.816
$ws.Range("M27").Value = -29079.816

$ws.Range("H46").Value = 1138.5264
$ws.Range("I46").Value = 622.2222
$ws.Range("J46").Value = 1603.2
$ws.Range("K46").Value = 1866.6666
$ws.Range("L46").Value = 4809.6
$ws.Range("M46").Value = -1775.6666
$ws.Range("N46").Value = -4991.6

$ws.Range("H64").Value = 2666.6667
$ws.Range("J64").Value = 3000
$ws.Range("L64").Value = 9000
$ws.Range("N64").Value = -9540

$ws.Range("H67").Value = 2666.6667
$ws.Range("J67").Value = 3000
$ws.Range("L67").Value = 9000
$ws.Range("N67").Value = -10872

$ws.Range("H98").Value = 697.5
$ws.Range("J98").Value = 711.4286
$ws.Range("L98").Value = 2134.2858
$ws.Range("N98").Value = -5130.2858

$ws.Range("H122").Value = 472
$ws.Range("I122").Value = 254
$ws.Range("K122").Value = 2286
$ws.Range("M122").Value = 164

$ws.Range("H129").Value = 1034.4166
$ws.Range("I129").Value = 849.5
$ws.Range("J129").Value = 1126.875
$ws.Range("K129").Value = 2548.5
$ws.Range("L129").Value = 3380.625
$ws.Range("M129").Value = 2451.5
$ws.Range("N129").Value = -13380.625

$ws.Range("H131").Value = 784.48
$ws.Range("J131").Value = 788.2782999999999
$ws.Range("L131").Value = 2364.8349
$ws.Range("N131").Value = -12444.8349

$ws.Range("H140").Value = 1586.1177
$ws.Range("I140").Value = 1368.6
$ws.Range("K140").Value = 4105.799999999999
$ws.Range("M140").Value = 1074.200000000001

# ----------------------------------------------------------------------
# LTW
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H16").Value = 357.2353
$ws.Range("I16").Value = 348.3125
$ws.Range("K16").Value = 348.3125
$ws.Range("M16").Value = -178.3125

$ws.Range("H22").Value = 5022.1113
$ws.Range("I22").Value = 6967
$ws.Range("J22").Value = 4049.6667
$ws.Range("K22").Value = 6967
$ws.Range("L22").Value = 4049.6667
$ws.Range("M22").Value = -6672
$ws.Range("N22").Value = -4639.6667

$ws.Range("H27").Value = 5022.1113
$ws.Range("I27").Value = 6967
$ws.Range("J27").Value = 4049.6667
$ws.Range("K27").Value = 6967
$ws.Range("L27").Value = 4049.6667
$ws.Range("M27").Value = -6860
$ws.Range("N27").Value = -4263.6667
